# Updates crypto price/volume figures to the latest scraped values.
# Rows 40/41 and 49/50 also had their coin entries swap ranking position
# (RenderToken/EthereumClassic and Mantle/InjectiveProtocol traded places).
#
# Column D ("Price") cells are forced to Text format before the write and
# reset back to the Normal style afterwards. Excel's COM layer otherwise
# reinterprets plain decimal-looking strings (e.g. "573.55") as floating
# point numbers, which corrupts the exact textual price formatting that the
# sheet relies on (and produces binary float noise like 573.54999999999995).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.762.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.459.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.458.85"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000176"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.907.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.697.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.460.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.53%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +18.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "650.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.593.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0973"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.74%  "
$ws.Range("E31").Value = "  +3.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("E34").Value = "  -4.45%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "152.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.368"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₆0314"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -70.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.38%  "
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.605"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0510"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
